$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028686040142017
$ws.Range("D2").Value = 1.033738532769151
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.027240847223345
$ws.Range("I2").Value = 1.035512461178826
$ws.Range("J2").Value = 1.033836983218924
$ws.Range("K2").Value = 1.036540017711374
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.03006115112178
$ws.Range("N2").Value = 1.035305150317866
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029595894531259
$ws.Range("D3").Value = 1.034435183756629
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.028782955642064
$ws.Range("I3").Value = 1.035756701744745
$ws.Range("J3").Value = 1.034387754605021
$ws.Range("K3").Value = 1.037046011267685
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.031408958097205
$ws.Range("N3").Value = 1.035856703862494
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030184224280929
$ws.Range("D4").Value = 1.034885490620382
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.02978039551877
$ws.Range("I4").Value = 1.03591310708797
$ws.Range("J4").Value = 1.034743118806055
$ws.Range("K4").Value = 1.037372277069489
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.032280191387182
$ws.Range("N4").Value = 1.036212572721454
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030431461001438
$ws.Range("D5").Value = 1.035074686105635
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.030199627043612
$ws.Range("I5").Value = 1.035978468562857
$ws.Range("J5").Value = 1.034892269547781
$ws.Range("K5").Value = 1.037509164905341
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.032646249850399
$ws.Range("N5").Value = 1.036361935274341
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030472967448976
$ws.Range("D6").Value = 1.035106446194165
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.030270012609348
$ws.Range("I6").Value = 1.035989420094244
$ws.Range("J6").Value = 1.034917298290334
$ws.Range("K6").Value = 1.037532132890088
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.032707700611583
$ws.Range("N6").Value = 1.03638699956058
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030187528251287
$ws.Range("D7").Value = 1.034888019106178
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.029785997666585
$ws.Range("I7").Value = 1.035913981988868
$ws.Range("J7").Value = 1.034745112725493
$ws.Range("K7").Value = 1.037374107249292
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.032285083489967
$ws.Range("N7").Value = 1.036214569472487
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028993614124524
$ws.Range("D8").Value = 1.033974066971093
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.027762097548601
$ws.Range("I8").Value = 1.035595341887312
$ws.Range("J8").Value = 1.034023330682774
$ws.Range("K8").Value = 1.036711257873264
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.030516835406754
$ws.Range("N8").Value = 1.035491762416477
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026886654202034
$ws.Range("D9").Value = 1.032359957131968
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.024192367797846
$ws.Range("I9").Value = 1.035021337245519
$ws.Range("J9").Value = 1.032743622191133
$ws.Range("K9").Value = 1.035534455099901
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.027393950896914
$ws.Range("N9").Value = 1.034210236591993
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025479885627131
$ws.Range("D10").Value = 1.031281471048425
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.021809953571761
$ws.Range("I10").Value = 1.034630246821288
$ws.Range("J10").Value = 1.0318851962094
$ws.Range("K10").Value = 1.034744013517697
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.025307029608644
$ws.Range("N10").Value = 1.033350591546903
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024870228056729
$ws.Range("D11").Value = 1.03081390435335
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.020777651311267
$ws.Range("I11").Value = 1.034458901847835
$ws.Range("J11").Value = 1.031512229386473
$ws.Range("K11").Value = 1.034400340707758
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.024402119440459
$ws.Range("N11").Value = 1.032977095068308
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024643695673127
$ws.Range("D12").Value = 1.030640142914478
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.020394096386697
$ws.Range("I12").Value = 1.034394955996377
$ws.Range("J12").Value = 1.031373502675721
$ws.Range("K12").Value = 1.034272473759659
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.024065800634514
$ws.Range("N12").Value = 1.032838171349715
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024692291197792
$ws.Range("D13").Value = 1.030677419239766
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.020476375389461
$ws.Range("I13").Value = 1.034408686208903
$ws.Range("J13").Value = 1.03140326864268
$ws.Range("K13").Value = 1.03429991123227
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.024137951069831
$ws.Range("N13").Value = 1.032867979587761
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024851504422826
$ws.Range("D14").Value = 1.030799542939859
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.020745948862362
$ws.Range("I14").Value = 1.034453622199073
$ws.Range("J14").Value = 1.031500766076883
$ws.Range("K14").Value = 1.034389775503391
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.024374323231323
$ws.Range("N14").Value = 1.032965615479504
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024949590534252
$ws.Range("D15").Value = 1.03087477594124
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.020912026951261
$ws.Range("I15").Value = 1.034481268894473
$ws.Range("J15").Value = 1.031560812241384
$ws.Range("K15").Value = 1.034445115804116
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.024519933928316
$ws.Range("N15").Value = 1.033025746916446
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025520335632339
$ws.Range("D16").Value = 1.03131248979812
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.021878448740983
$ws.Range("I16").Value = 1.034641576264396
$ws.Range("J16").Value = 1.03190992214125
$ws.Range("K16").Value = 1.03476679229639
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.02536705849839
$ws.Range("N16").Value = 1.033375352592413
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025878209937022
$ws.Range("D17").Value = 1.031586902206338
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.022484467254391
$ws.Range("I17").Value = 1.034741597192314
$ws.Range("J17").Value = 1.032128571132608
$ws.Range("K17").Value = 1.034968194743793
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.025898096180087
$ws.Range("N17").Value = 1.033594312090418
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026086902079054
$ws.Range("D18").Value = 1.031746906829353
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.02283788038826
$ws.Range("I18").Value = 1.034799744726093
$ws.Range("J18").Value = 1.03225598359238
$ws.Range("K18").Value = 1.035085533706508
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.026207720334865
$ws.Range("N18").Value = 1.033721905490497
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026158052265194
$ws.Range("D19").Value = 1.031801454856693
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.022958373932563
$ws.Range("I19").Value = 1.034819538805041
$ws.Range("J19").Value = 1.032299407305158
$ws.Range("K19").Value = 1.035125520255312
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.026313273794328
$ws.Range("N19").Value = 1.033765390869928
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025839818575502
$ws.Range("D20").Value = 1.031557466091216
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.022419454201887
$ws.Range("I20").Value = 1.03473088585554
$ws.Range("J20").Value = 1.032105124766795
$ws.Range("K20").Value = 1.034946600202799
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.025841133432865
$ws.Range("N20").Value = 1.033570832428077
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024804622266462
$ws.Range("D21").Value = 1.030763582937354
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.020666569335401
$ws.Range("I21").Value = 1.034440397972259
$ws.Range("J21").Value = 1.031472060767967
$ws.Range("K21").Value = 1.034363318557792
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.02430472294635
$ws.Range("N21").Value = 1.032936869405756
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024153299224367
$ws.Range("D22").Value = 1.03026393751225
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.019563810347733
$ws.Range("I22").Value = 1.034256016861596
$ws.Range("J22").Value = 1.031072927450235
$ws.Range("K22").Value = 1.033995362049142
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.023337590365613
$ws.Range("N22").Value = 1.032537169272932
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024498621156335
$ws.Range("D23").Value = 1.030528856345645
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.020148467630061
$ws.Range("I23").Value = 1.034353925758226
$ws.Range("J23").Value = 1.03128462009824
$ws.Range("K23").Value = 1.034190538888332
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.023850394684434
$ws.Range("N23").Value = 1.032749162548779
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025857166133091
$ws.Range("D24").Value = 1.031570767175906
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.022448831007998
$ws.Range("I24").Value = 1.034735726443722
$ws.Range("J24").Value = 1.03211571954638
$ws.Range("K24").Value = 1.034956358265394
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.025866872819355
$ws.Range("N24").Value = 1.033581442253464
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027431727250192
$ws.Range("D25").Value = 1.032777669467803
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.025115661071168
$ws.Range("I25").Value = 1.035171215170363
$ws.Range("J25").Value = 1.03307538799971
$ws.Range("K25").Value = 1.035839727469045
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.028202151170931
$ws.Range("N25").Value = 1.034542473546073
